$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# orders.xlsx: row 46 MI column had a floating-point -0 artifact; normalize to 0
$ws.Range("MI46").Value = 0

# New row 47 (2020-01-27 trading day). Give A47 the same bold/border/centered
# "YYYY-MM-DD HH:MM:SS" date style used by the rest of column A (style index 2)
# by copying formats from the cell directly above it, then fill in the values.
$ws.Range("A46").Copy()
$ws.Range("A47").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate new row 47 with data (one cell per ticker/return column, A:RR)
$ws.Range("A47").Value = 43857
$ws.Range("B47").Value = 0
$ws.Range("C47").Value = 1.268419838764061
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("F47").Value = 2.358713936863353
$ws.Range("G47").Value = 0
$ws.Range("H47").Value = 1.638496125706936
$ws.Range("I47").Value = 25.1261541894969
$ws.Range("J47").Value = 0.9165119296602739
$ws.Range("K47").Value = 13.29371723668544
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = 0
$ws.Range("N47").Value = 0.1464025497482977
$ws.Range("O47").Value = 8.20246104638511
$ws.Range("P47").Value = 17.07374710982651
$ws.Range("Q47").Value = 16.64710303647132
$ws.Range("R47").Value = 0.1436804500462117
$ws.Range("S47").Value = 0
$ws.Range("T47").Value = 0
$ws.Range("U47").Value = 0.9850927460414027
$ws.Range("V47").Value = 0
$ws.Range("W47").Value = 0
$ws.Range("X47").Value = 1.616170784843689
$ws.Range("Y47").Value = 3.304059836047912
$ws.Range("Z47").Value = 0.8517291555195001
$ws.Range("AA47").Value = 0.8579526519478549
$ws.Range("AB47").Value = 0
$ws.Range("AC47").Value = 0
$ws.Range("AD47").Value = -107.0858776254204
$ws.Range("AE47").Value = 1.144110600044883
$ws.Range("AF47").Value = 98.43219405584477
$ws.Range("AG47").Value = 2.080364991291731
$ws.Range("AH47").Value = 9.406551601124818
$ws.Range("AI47").Value = 0
$ws.Range("AJ47").Value = 0.008579564785499372
$ws.Range("AK47").Value = 1.641211167097538
$ws.Range("AL47").Value = 9.857679230287033
$ws.Range("AM47").Value = 0
$ws.Range("AN47").Value = 0
$ws.Range("AO47").Value = 0
$ws.Range("AP47").Value = 0
$ws.Range("AQ47").Value = 22.0343803933103
$ws.Range("AR47").Value = 4.655858679254465
$ws.Range("AS47").Value = 0
$ws.Range("AT47").Value = 1.819952300264561
$ws.Range("AU47").Value = 0
$ws.Range("AV47").Value = 0
$ws.Range("AW47").Value = 3.395836583427013
$ws.Range("AX47").Value = 5.442271246497171
$ws.Range("AY47").Value = 34.81575238964069
$ws.Range("AZ47").Value = 17.40425990141739
$ws.Range("BA47").Value = 0
$ws.Range("BB47").Value = -3.163679351483246
$ws.Range("BC47").Value = 2.765393963994057
$ws.Range("BD47").Value = 3.414478491757052
$ws.Range("BE47").Value = 0
$ws.Range("BF47").Value = 3.665160490380629
$ws.Range("BG47").Value = 0.393946620058621
$ws.Range("BH47").Value = 2.3510078273749
$ws.Range("BI47").Value = 1.45590942009224
$ws.Range("BJ47").Value = 1.050593818650505
$ws.Range("BK47").Value = 0
$ws.Range("BL47").Value = 1.65302361250707
$ws.Range("BM47").Value = 54.90287838440236
$ws.Range("BN47").Value = -16.61381925252601
$ws.Range("BO47").Value = 0
$ws.Range("BP47").Value = 15.00772204252439
$ws.Range("BQ47").Value = 0.1877693186332436
$ws.Range("BR47").Value = -0.763684877859248
$ws.Range("BS47").Value = 2.999449431736451
$ws.Range("BT47").Value = 0
$ws.Range("BU47").Value = 0
$ws.Range("BV47").Value = 2.632176162880043
$ws.Range("BW47").Value = 0.9186398063713455
$ws.Range("BX47").Value = 0.05750638838344679
$ws.Range("BY47").Value = 1.508084318399497
$ws.Range("BZ47").Value = 1.579031101077192
$ws.Range("CA47").Value = 0
$ws.Range("CB47").Value = 0
$ws.Range("CC47").Value = 8.532036254643543
$ws.Range("CD47").Value = 0
$ws.Range("CE47").Value = 12.81566084858719
$ws.Range("CF47").Value = 0
$ws.Range("CG47").Value = 0
$ws.Range("CH47").Value = 19.17224926885285
$ws.Range("CI47").Value = 0.3722449905637859
$ws.Range("CJ47").Value = 10.48636737969889
$ws.Range("CK47").Value = 0
$ws.Range("CL47").Value = 3.682434396701979
$ws.Range("CM47").Value = 4.507505136475174
$ws.Range("CN47").Value = 0
$ws.Range("CO47").Value = 0
$ws.Range("CP47").Value = 0.3173573337822404
$ws.Range("CQ47").Value = 0
$ws.Range("CR47").Value = 11.19095550405962
$ws.Range("CS47").Value = 0
$ws.Range("CT47").Value = 43.50376528829895
$ws.Range("CU47").Value = 0
$ws.Range("CV47").Value = 6.180993899119358
$ws.Range("CW47").Value = 0
$ws.Range("CX47").Value = 0
$ws.Range("CY47").Value = 0
$ws.Range("CZ47").Value = 0.1232755303126538
$ws.Range("DA47").Value = 0
$ws.Range("DB47").Value = 0.2881486938136106
$ws.Range("DC47").Value = 0
$ws.Range("DD47").Value = 0
$ws.Range("DE47").Value = 0
$ws.Range("DF47").Value = 0
$ws.Range("DG47").Value = 0
$ws.Range("DH47").Value = 2.956335600175862
$ws.Range("DI47").Value = 18.4988975339503
$ws.Range("DJ47").Value = 0
$ws.Range("DK47").Value = 0
$ws.Range("DL47").Value = 4.070164642470672
$ws.Range("DM47").Value = 0
$ws.Range("DN47").Value = -4.831078302068477
$ws.Range("DO47").Value = 2.98688495368981
$ws.Range("DP47").Value = 3.351180454461883
$ws.Range("DQ47").Value = 4.752961451060798
$ws.Range("DR47").Value = 0
$ws.Range("DS47").Value = 0
$ws.Range("DT47").Value = 0
$ws.Range("DU47").Value = 0
$ws.Range("DV47").Value = 0
$ws.Range("DW47").Value = 0
$ws.Range("DX47").Value = 1.150724968207783
$ws.Range("DY47").Value = 1.037238854404848
$ws.Range("DZ47").Value = 0
$ws.Range("EA47").Value = 4.072231449632454
$ws.Range("EB47").Value = 0
$ws.Range("EC47").Value = 59.82858333979311
$ws.Range("ED47").Value = 0
$ws.Range("EE47").Value = 0
$ws.Range("EF47").Value = -0.9738272266295667
$ws.Range("EG47").Value = 10.09057669165986
$ws.Range("EH47").Value = 1.28735113369499
$ws.Range("EI47").Value = 0
$ws.Range("EJ47").Value = -17.36231960163366
$ws.Range("EK47").Value = 0
$ws.Range("EL47").Value = 2.559549231195433
$ws.Range("EM47").Value = 6.736624703175551
$ws.Range("EN47").Value = 0
$ws.Range("EO47").Value = 0
$ws.Range("EP47").Value = 0
$ws.Range("EQ47").Value = 0
$ws.Range("ER47").Value = 0
$ws.Range("ES47").Value = -2.610058757624586
$ws.Range("ET47").Value = 22.23540089072844
$ws.Range("EU47").Value = 0
$ws.Range("EV47").Value = 0
$ws.Range("EW47").Value = -0.320013916568314
$ws.Range("EX47").Value = 0
$ws.Range("EY47").Value = 0
$ws.Range("EZ47").Value = 0.05135207115939711
$ws.Range("FA47").Value = 4.679817816659863
$ws.Range("FB47").Value = 0
$ws.Range("FC47").Value = 0
$ws.Range("FD47").Value = 1.106913875240934
$ws.Range("FE47").Value = 0.01032414515469871
$ws.Range("FF47").Value = 0
$ws.Range("FG47").Value = 46.16880139988689
$ws.Range("FH47").Value = 0
$ws.Range("FI47").Value = 0
$ws.Range("FJ47").Value = 12.51138141524359
$ws.Range("FK47").Value = 0
$ws.Range("FL47").Value = 1.779056133720587
$ws.Range("FM47").Value = 0
$ws.Range("FN47").Value = 0
$ws.Range("FO47").Value = 8.014525246314463
$ws.Range("FP47").Value = 10.59612953715202
$ws.Range("FQ47").Value = 3.59570388584865
$ws.Range("FR47").Value = 0
$ws.Range("FS47").Value = 2.304506186047945
$ws.Range("FT47").Value = -0.3059452754749401
$ws.Range("FU47").Value = 0
$ws.Range("FV47").Value = 0.04517168281969219
$ws.Range("FW47").Value = 0
$ws.Range("FX47").Value = -4.032532155926219
$ws.Range("FY47").Value = 2.499588163234762
$ws.Range("FZ47").Value = 2.437686733436465
$ws.Range("GA47").Value = 10.60691622907962
$ws.Range("GB47").Value = 10.38445306631064
$ws.Range("GC47").Value = 0
$ws.Range("GD47").Value = 0
$ws.Range("GE47").Value = 0
$ws.Range("GF47").Value = 0
$ws.Range("GG47").Value = 0
$ws.Range("GH47").Value = 12.35320994870676
$ws.Range("GI47").Value = 1.330779943590386
$ws.Range("GJ47").Value = 0
$ws.Range("GK47").Value = 3.868687787732938
$ws.Range("GL47").Value = 18.92079399096087
$ws.Range("GM47").Value = 4.745196288303987
$ws.Range("GN47").Value = 0
$ws.Range("GO47").Value = 5.121498852822072
$ws.Range("GP47").Value = 2.337966704698957
$ws.Range("GQ47").Value = 0.1124789352095243
$ws.Range("GR47").Value = 0
$ws.Range("GS47").Value = 0
$ws.Range("GT47").Value = 0
$ws.Range("GU47").Value = 12.80132163098972
$ws.Range("GV47").Value = 0
$ws.Range("GW47").Value = 21.14173552325769
$ws.Range("GX47").Value = 0
$ws.Range("GY47").Value = 0
$ws.Range("GZ47").Value = 13.53522038860433
$ws.Range("HA47").Value = 0
$ws.Range("HB47").Value = 0
$ws.Range("HC47").Value = 0
$ws.Range("HD47").Value = 0.02447228491890208
$ws.Range("HE47").Value = 2.748281111006264
$ws.Range("HF47").Value = 8.225881438869692
$ws.Range("HG47").Value = 0
$ws.Range("HH47").Value = 0
$ws.Range("HI47").Value = 0
$ws.Range("HJ47").Value = -9.562638085978961
$ws.Range("HK47").Value = 4.955718679490701
$ws.Range("HL47").Value = 2.667945642440429
$ws.Range("HM47").Value = 0
$ws.Range("HN47").Value = 0
$ws.Range("HO47").Value = 1.238194394947314
$ws.Range("HP47").Value = 0
$ws.Range("HQ47").Value = 60.42452968250382
$ws.Range("HR47").Value = 13.72900862817323
$ws.Range("HS47").Value = 0
$ws.Range("HT47").Value = 2.687140375700949
$ws.Range("HU47").Value = 0
$ws.Range("HV47").Value = 2.115556159813877
$ws.Range("HW47").Value = 0
$ws.Range("HX47").Value = 4.1450853793433
$ws.Range("HY47").Value = 2.168724104753551
$ws.Range("HZ47").Value = 0
$ws.Range("IA47").Value = 64.51524708733905
$ws.Range("IB47").Value = 17.24777301880943
$ws.Range("IC47").Value = 0
$ws.Range("ID47").Value = 0
$ws.Range("IE47").Value = 0.6061316550684808
$ws.Range("IF47").Value = 0
$ws.Range("IG47").Value = 0.07951605170156029
$ws.Range("IH47").Value = 0
$ws.Range("II47").Value = 28.91381123878296
$ws.Range("IJ47").Value = 0
$ws.Range("IK47").Value = 0
$ws.Range("IL47").Value = 0
$ws.Range("IM47").Value = 0
$ws.Range("IN47").Value = 0.3653443083019106
$ws.Range("IO47").Value = 2.786555496449182
$ws.Range("IP47").Value = 3.283650231119935
$ws.Range("IQ47").Value = 0
$ws.Range("IR47").Value = 0
$ws.Range("IS47").Value = 32.08437532294272
$ws.Range("IT47").Value = -2.443409136816172
$ws.Range("IU47").Value = 3.182789288040681
$ws.Range("IV47").Value = 0
$ws.Range("IW47").Value = 16.87311744421061
$ws.Range("IX47").Value = 0
$ws.Range("IY47").Value = 4.806228152620179
$ws.Range("IZ47").Value = 1.241850367911979
$ws.Range("JA47").Value = 3.638988104308069
$ws.Range("JB47").Value = 48.11003940757541
$ws.Range("JC47").Value = 28.06494571558505
$ws.Range("JD47").Value = 0.04000452812339095
$ws.Range("JE47").Value = -4.557673389249885
$ws.Range("JF47").Value = 4.281414533268844
$ws.Range("JG47").Value = 8.420200852132496
$ws.Range("JH47").Value = 0.341467418137583
$ws.Range("JI47").Value = 0.4134063288629477
$ws.Range("JJ47").Value = 0
$ws.Range("JK47").Value = 36.26321131768464
$ws.Range("JL47").Value = 0
$ws.Range("JM47").Value = 0
$ws.Range("JN47").Value = 2.228084350133997
$ws.Range("JO47").Value = 0
$ws.Range("JP47").Value = 69.6092845210419
$ws.Range("JQ47").Value = 0
$ws.Range("JR47").Value = -0.07317037373683277
$ws.Range("JS47").Value = 3.99224225405419
$ws.Range("JT47").Value = 26.33325583540818
$ws.Range("JU47").Value = 16.65941073292936
$ws.Range("JV47").Value = 3.16277070465938
$ws.Range("JW47").Value = 0
$ws.Range("JX47").Value = -1.84243651009848
$ws.Range("JY47").Value = 35.52772997703073
$ws.Range("JZ47").Value = 0.4463474930100233
$ws.Range("KA47").Value = -0.243388704379953
$ws.Range("KB47").Value = 7.698266383481382
$ws.Range("KC47").Value = 0
$ws.Range("KD47").Value = 0
$ws.Range("KE47").Value = 8.656140379831413
$ws.Range("KF47").Value = 0
$ws.Range("KG47").Value = -8.126878119515482
$ws.Range("KH47").Value = 0.4839390822741052
$ws.Range("KI47").Value = 0
$ws.Range("KJ47").Value = 0
$ws.Range("KK47").Value = 0
$ws.Range("KL47").Value = 2.085156638345111
$ws.Range("KM47").Value = 0.01326316314913356
$ws.Range("KN47").Value = 0
$ws.Range("KO47").Value = 1.874128976166162
$ws.Range("KP47").Value = 0
$ws.Range("KQ47").Value = 5.679945798114034
$ws.Range("KR47").Value = 20.96033614063731
$ws.Range("KS47").Value = 0
$ws.Range("KT47").Value = 0
$ws.Range("KU47").Value = 0
$ws.Range("KV47").Value = 4.182197204893356
$ws.Range("KW47").Value = 0.7346426906021009
$ws.Range("KX47").Value = 0.004962519801303023
$ws.Range("KY47").Value = 7.579631998538161
$ws.Range("KZ47").Value = 4.030721054415721
$ws.Range("LA47").Value = -0.06632527732671178
$ws.Range("LB47").Value = 33.0115843254581
$ws.Range("LC47").Value = 2.082527446535693
$ws.Range("LD47").Value = 0
$ws.Range("LE47").Value = 0.07480146224122564
$ws.Range("LF47").Value = 0
$ws.Range("LG47").Value = 2.03623736702832
$ws.Range("LH47").Value = 0
$ws.Range("LI47").Value = 0
$ws.Range("LJ47").Value = 18.20167371440328
$ws.Range("LK47").Value = 35.78497493219777
$ws.Range("LL47").Value = 0.2493163912023242
$ws.Range("LM47").Value = 0.329735507030108
$ws.Range("LN47").Value = 4.055732498288705
$ws.Range("LO47").Value = 2.272273380858223
$ws.Range("LP47").Value = 0
$ws.Range("LQ47").Value = 6.196571735125872
$ws.Range("LR47").Value = 0.1878670710616888
$ws.Range("LS47").Value = 0
$ws.Range("LT47").Value = 0
$ws.Range("LU47").Value = 0.03934841065151673
$ws.Range("LV47").Value = 1.151492524779343
$ws.Range("LW47").Value = 81.01749659023199
$ws.Range("LX47").Value = 1.11823753611867
$ws.Range("LY47").Value = 0
$ws.Range("LZ47").Value = 0.225221430241926
$ws.Range("MA47").Value = 1.159043403789205
$ws.Range("MB47").Value = 41.81912345837327
$ws.Range("MC47").Value = 0
$ws.Range("MD47").Value = 0
$ws.Range("ME47").Value = 0
$ws.Range("MF47").Value = -1.959306097516048
$ws.Range("MG47").Value = 50.43125979874253
$ws.Range("MH47").Value = 0
$ws.Range("MI47").Value = -0
$ws.Range("MJ47").Value = 0
$ws.Range("MK47").Value = 0
$ws.Range("ML47").Value = 3.535739695869722
$ws.Range("MM47").Value = 0
$ws.Range("MN47").Value = 0.3035007166577799
$ws.Range("MO47").Value = 47.76981646536183
$ws.Range("MP47").Value = 6.686123320477719
$ws.Range("MQ47").Value = 46.87450368390637
$ws.Range("MR47").Value = 0
$ws.Range("MS47").Value = 9.777562093076767
$ws.Range("MT47").Value = 0.1238730175267051
$ws.Range("MU47").Value = 0
$ws.Range("MV47").Value = 0
$ws.Range("MW47").Value = 0
$ws.Range("MX47").Value = 0
$ws.Range("MY47").Value = 2.641189323312148
$ws.Range("MZ47").Value = 0
$ws.Range("NA47").Value = 0.109246662570623
$ws.Range("NB47").Value = 7.621615506374667
$ws.Range("NC47").Value = 0
$ws.Range("ND47").Value = 6.980425572958097
$ws.Range("NE47").Value = 3.68577200146342
$ws.Range("NF47").Value = 1.140846576378586
$ws.Range("NG47").Value = 24.29338005741101
$ws.Range("NH47").Value = 1.902729005427659
$ws.Range("NI47").Value = 0
$ws.Range("NJ47").Value = 0
$ws.Range("NK47").Value = -5.852162141685881
$ws.Range("NL47").Value = 0
$ws.Range("NM47").Value = 0
$ws.Range("NN47").Value = -2.808898812117377
$ws.Range("NO47").Value = 0
$ws.Range("NP47").Value = 4.735716730866955
$ws.Range("NQ47").Value = 9.916010074399821
$ws.Range("NR47").Value = 13.28426046306231
$ws.Range("NS47").Value = 0
$ws.Range("NT47").Value = 0
$ws.Range("NU47").Value = 0
$ws.Range("NV47").Value = 0.03213895108101061
$ws.Range("NW47").Value = 0.1331444284576846
$ws.Range("NX47").Value = 18.56109705514234
$ws.Range("NY47").Value = 0.01373058642904645
$ws.Range("NZ47").Value = -16.2973955922464
$ws.Range("OA47").Value = 0
$ws.Range("OB47").Value = 15.4718271500322
$ws.Range("OC47").Value = 0
$ws.Range("OD47").Value = 0
$ws.Range("OE47").Value = -2.401042571934994
$ws.Range("OF47").Value = 0
$ws.Range("OG47").Value = 0
$ws.Range("OH47").Value = 2.929996818911064
$ws.Range("OI47").Value = 0
$ws.Range("OJ47").Value = 1.57070369563877
$ws.Range("OK47").Value = 0
$ws.Range("OL47").Value = 0
$ws.Range("OM47").Value = 0
$ws.Range("ON47").Value = 0
$ws.Range("OO47").Value = 0
$ws.Range("OP47").Value = 0.2306361317415169
$ws.Range("OQ47").Value = 14.95476814511221
$ws.Range("OR47").Value = 0.007333998321566781
$ws.Range("OS47").Value = 1.0355074853139
$ws.Range("OT47").Value = 0
$ws.Range("OU47").Value = 0.6607557905098105
$ws.Range("OV47").Value = 0
$ws.Range("OW47").Value = 0.01907174127433109
$ws.Range("OX47").Value = 0
$ws.Range("OY47").Value = 32.15519289984081
$ws.Range("OZ47").Value = 0
$ws.Range("PA47").Value = 1.687445795827358
$ws.Range("PB47").Value = 0
$ws.Range("PC47").Value = 1.265070704930949
$ws.Range("PD47").Value = 0
$ws.Range("PE47").Value = 20.2065282182129
$ws.Range("PF47").Value = 4.210174945785525
$ws.Range("PG47").Value = 0
$ws.Range("PH47").Value = -3.218573468129193
$ws.Range("PI47").Value = 3.943715871801714
$ws.Range("PJ47").Value = 0
$ws.Range("PK47").Value = 17.46141349550214
$ws.Range("PL47").Value = 0
$ws.Range("PM47").Value = 0
$ws.Range("PN47").Value = -0.1895550809639275
$ws.Range("PO47").Value = 8.638311860984345
$ws.Range("PP47").Value = 10.59887978475626
$ws.Range("PQ47").Value = 6.013101974591819
$ws.Range("PR47").Value = 1.403186607079306
$ws.Range("PS47").Value = 0.8061222668388552
$ws.Range("PT47").Value = 8.54881624014439
$ws.Range("PU47").Value = 0
$ws.Range("PV47").Value = 0
$ws.Range("PW47").Value = -4.321856666234282
$ws.Range("PX47").Value = 3.481372447755007
$ws.Range("PY47").Value = 0
$ws.Range("PZ47").Value = 0.07820352214579529
$ws.Range("QA47").Value = 13.7303321948902
$ws.Range("QB47").Value = 0
$ws.Range("QC47").Value = 0
$ws.Range("QD47").Value = 0
$ws.Range("QE47").Value = 3.406877505814236
$ws.Range("QF47").Value = -0.001355979304939636
$ws.Range("QG47").Value = 0.6388854716640679
$ws.Range("QH47").Value = 0
$ws.Range("QI47").Value = 14.02384230265102
$ws.Range("QJ47").Value = 0
$ws.Range("QK47").Value = 13.36342878677908
$ws.Range("QL47").Value = 0.6877764559268087
$ws.Range("QM47").Value = -0.004525757986332923
$ws.Range("QN47").Value = 3.918954110351478
$ws.Range("QO47").Value = 43.12969935699857
$ws.Range("QP47").Value = 25.72061733965666
$ws.Range("QQ47").Value = 6.58282300356376
$ws.Range("QR47").Value = 0
$ws.Range("QS47").Value = 0
$ws.Range("QT47").Value = 0
$ws.Range("QU47").Value = 4.53177374167079
$ws.Range("QV47").Value = 0
$ws.Range("QW47").Value = 5.537733889190633
$ws.Range("QX47").Value = 0
$ws.Range("QY47").Value = 0
$ws.Range("QZ47").Value = 0
$ws.Range("RA47").Value = 4.838480586284447
$ws.Range("RB47").Value = 0
$ws.Range("RC47").Value = 1.466331902930548
$ws.Range("RD47").Value = 2.830066585283987
$ws.Range("RE47").Value = 5.105724601786505
$ws.Range("RF47").Value = 0
$ws.Range("RG47").Value = 0
$ws.Range("RH47").Value = 0
$ws.Range("RI47").Value = 8.991837814357496
$ws.Range("RJ47").Value = 3.262147661126562
$ws.Range("RK47").Value = 4.384218087098361
$ws.Range("RL47").Value = 20.80010295201873
$ws.Range("RM47").Value = -1.197928908976706
$ws.Range("RN47").Value = 0
$ws.Range("RO47").Value = 0
$ws.Range("RP47").Value = -0.4917433900604919
$ws.Range("RQ47").Value = 0
$ws.Range("RR47").Value = 2.711827948810082
